$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two shared-string values
$ws.Range("A1").Value = "surabhisanjan05"
$ws.Range("B1").Value = "Theend@1"

# Update the active cell selection on the sheet
$ws.Range("B1").Select()
